$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row ---
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Clear old columns F and G (no longer used)
$ws.Range("F1:G4").Clear()

# --- Data ---
$data = @(
    @(0, 3435.566666666667, 3512, 3399, 0.03216848373413086),
    @(1, 3234.333333333333, 3387, 3029, 0.03517893155415853),
    @(2, 3458.2,            3585, 3238, 0.03489200274149577),
    @(3, 3570.733333333333, 3917, 3312, 0.03482209841410319),
    @(4, 3298,              3555, 2981, 0.03622381687164307),
    @(5, 4081.466666666667, 4227, 3990, 0.04528450171152751),
    @(6, 3107.2,            3344, 2966, 0.04557886123657227),
    @(7, 3790.366666666667, 4021, 3582, 0.04214363098144532),
    @(8, 3149.266666666667, 3377, 2935, 0.04198220570882161),
    @(9, 3122.9,            3531, 3029, 0.03792573610941569)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}
